# Updated cryptos list on Thu Oct 19 19:30:23 UTC 2023 with GitHub Actions
#
# NOTE: a handful of "Price" values (column D) are single-decimal numeric
# strings (e.g. "210.33") that Excel's type-inference would otherwise coerce
# into floating-point numbers (losing exact text formatting / precision).
# Those are written with a leading apostrophe so Excel keeps them as text,
# matching the source data exactly. Multi-dot / non-numeric price strings
# (e.g. "28.655.85") need no such treatment since Excel can't parse them as
# numbers anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.655.85'
$ws.Range("E2").Value = '  +1.08%  '
$ws.Range("D3").Value = '1.563.09'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").Value = "'210.33"
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("E7").Value = '  -0.48%  '
$ws.Range("D8").Value = "'25.10"
$ws.Range("E8").Value = '  +5.41%  '
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("E10").Value = '  -0.24%  '
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("D12").Value = '1.785.73'
$ws.Range("E12").Value = '  -0.48%  '
$ws.Range("D13").Value = '1.556.05'
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").Value = '28.669.97'
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("D15").Value = "'0.516"
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("E16").Value = '  -0.81%  '
$ws.Range("D17").Value = "'61.27"
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = "'228.45"
$ws.Range("E18").Value = '  +0.39%  '
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("D20").Value = '0.0₃0678'
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D22").Value = "'3.91"
$ws.Range("E22").Value = '  -0.59%  '
$ws.Range("D23").Value = "'9.02"
$ws.Range("E23").Value = '  +0.87%  '
$ws.Range("D24").Value = "'2.07"
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("D25").Value = "'151.26"
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("D26").Value = "'14.75"
$ws.Range("E26").Value = '  -0.97%  '
$ws.Range("E27").Value = '  +0.54%  '
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("E29").Value = '  -1.75%  '
$ws.Range("E30").Value = '  -3.95%  '
$ws.Range("E31").Value = '  -2.66%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").Value = '1.392.27'
$ws.Range("E33").Value = '  +1.03%  '
$ws.Range("E34").Value = '  -3.24%  '
$ws.Range("D35").Value = "'1.03"
$ws.Range("E35").Value = '  -4.07%  '
$ws.Range("E36").Value = '  -1.39%  '
$ws.Range("E37").Value = '  +1.92%  '
$ws.Range("E38").Value = '  -2.25%  '
$ws.Range("D39").Value = "'0.0161"
$ws.Range("E39").Value = '  -0.82%  '
$ws.Range("E40").Value = '  +1.83%  '
$ws.Range("D41").Value = "'0.518"
$ws.Range("E41").Value = '  -0.40%  '
$ws.Range("E42").Value = '  -0.38%  '
$ws.Range("D43").Value = "'0.771"
$ws.Range("E43").Value = '  -1.49%  '
$ws.Range("E44").Value = '  -2.90%  '
$ws.Range("D45").Value = "'64.03"
$ws.Range("E45").Value = '  +2.93%  '
$ws.Range("E46").Value = '  -1.81%  '
$ws.Range("D47").Value = '1.697.95'
$ws.Range("E47").Value = '  -0.53%  '
$ws.Range("D48").Value = "'0.871"
$ws.Range("E48").Value = '  -5.08%  '
$ws.Range("D49").Value = "'85.11"
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("D50").Value = "'43.23"
$ws.Range("E50").Value = '  +7.29%  '
$ws.Range("E51").Value = '  +1.01%  '
